$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (GoalPoseY) value changed from 255 to 208
$ws.Range("B12").Value = 208

# Active selection moved from B14 to B13
$ws.Range("B13").Select()
